# Adds four new "Summary_*" worksheets (Summary_42..Summary_45) at the end
# of the workbook, continuing the existing repeating 4-sheet cycle of
# Summary tabs (full TestData summary / header-only / TestData summary with
# extra row / Summary_1-based summary).

$wb = $excel.ActiveWorkbook

function Add-SheetAtEnd {
    param([string]$name)
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $name
    return $newSheet
}

# ---------------------------------------------------------------------------
# Summary_42  (mirrors Summary_38 / Summary_2 pattern: full TestData summary)
# ---------------------------------------------------------------------------
$ws = Add-SheetAtEnd "Summary_42"

$ws.Range("A1").Value = "Column Heading"
$ws.Range("B1").Value = "Count"
$ws.Range("C1").Value = "Total"
$ws.Range("D1").Value = "Percentage"

$ws.Range("A2").Value = "Name"
$ws.Range("B2").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("C2").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("D2").Formula = "=B2/C2"
$ws.Range("D2").NumberFormat = "0.00%"

$ws.Range("A3").Value = "Age"
$ws.Range("B3").Formula = "=COUNTA(TestData!B:B)-1"
$ws.Range("C3").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("D3").Formula = "=B3/C3"
$ws.Range("D3").NumberFormat = "0.00%"

$ws.Range("A4").Value = "City"
$ws.Range("B4").Formula = "=COUNTA(TestData!C:C)-1"
$ws.Range("C4").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("D4").Formula = "=B4/C4"
$ws.Range("D4").NumberFormat = "0.00%"

$ws.Range("A5").Value = "Score"
$ws.Range("B5").Formula = "=COUNTA(TestData!D:D)-1"
$ws.Range("C5").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("D5").Formula = "=B5/C5"
$ws.Range("D5").NumberFormat = "0.00%"

$ws.Range("A6").Value = "Comments"
$ws.Range("B6").Formula = "=COUNTA(TestData!E:E)-1"
$ws.Range("C6").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("D6").Formula = "=B6/C6"
$ws.Range("D6").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# Summary_43  (mirrors Summary_39 / Summary_7 pattern: header row only)
# ---------------------------------------------------------------------------
$ws = Add-SheetAtEnd "Summary_43"

$ws.Range("A1").Value = "Column Heading"
$ws.Range("B1").Value = "Count"
$ws.Range("C1").Value = "Total"
$ws.Range("D1").Value = "Percentage"

# ---------------------------------------------------------------------------
# Summary_44  (mirrors Summary_40 / Summary_8 pattern: TestData summary + extra row)
# ---------------------------------------------------------------------------
$ws = Add-SheetAtEnd "Summary_44"

$ws.Range("A1").Value = "Column Heading"
$ws.Range("B1").Value = "Count"
$ws.Range("C1").Value = "Total"
$ws.Range("D1").Value = "Percentage"

$ws.Range("A2").Value = "Name"
$ws.Range("B2").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("C2").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("D2").Formula = "=B2/C2"
$ws.Range("D2").NumberFormat = "0.00%"

$ws.Range("A3").Value = "Age"
$ws.Range("B3").Formula = "=COUNTA(TestData!B:B)-1"
$ws.Range("C3").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("D3").Formula = "=B3/C3"
$ws.Range("D3").NumberFormat = "0.00%"

$ws.Range("A4").Value = "City"
$ws.Range("B4").Formula = "=COUNTA(TestData!C:C)-1"
$ws.Range("C4").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("D4").Formula = "=B4/C4"
$ws.Range("D4").NumberFormat = "0.00%"

$ws.Range("A5").Value = "Score"
$ws.Range("B5").Formula = "=COUNTA(TestData!D:D)-1"
$ws.Range("C5").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("D5").Formula = "=B5/C5"
$ws.Range("D5").NumberFormat = "0.00%"

$ws.Range("A6").Value = "Comments"
$ws.Range("B6").Formula = "=COUNTA(TestData!E:E)-1"
$ws.Range("C6").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("D6").Formula = "=B6/C6"
$ws.Range("D6").NumberFormat = "0.00%"

# Row 7 has no label in column A.
$ws.Range("B7").Formula = "=COUNTA(TestData!F:F)-1"
$ws.Range("C7").Formula = "=COUNTA(TestData!A:A)-1"
$ws.Range("D7").Formula = "=B7/C7"
$ws.Range("D7").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# Summary_45  (mirrors Summary_41 / Summary_9 pattern: summary of Summary_1)
# ---------------------------------------------------------------------------
$ws = Add-SheetAtEnd "Summary_45"

$ws.Range("A1").Value = "Column Heading"
$ws.Range("B1").Value = "Count"
$ws.Range("C1").Value = "Total"
$ws.Range("D1").Value = "Percentage"

$ws.Range("A2").Value = "Column Heading"
$ws.Range("B2").Formula = "=COUNTA(Summary_1!A:A)-1"
$ws.Range("C2").Formula = "=COUNTA(Summary_1!A:A)-1"
$ws.Range("D2").Formula = "=B2/C2"
$ws.Range("D2").NumberFormat = "0.00%"

$ws.Range("A3").Value = "Count"
$ws.Range("B3").Formula = "=COUNTA(Summary_1!B:B)-1"
$ws.Range("C3").Formula = "=COUNTA(Summary_1!A:A)-1"
$ws.Range("D3").Formula = "=B3/C3"
$ws.Range("D3").NumberFormat = "0.00%"

$ws.Range("A4").Value = "Total"
$ws.Range("B4").Formula = "=COUNTA(Summary_1!C:C)-1"
$ws.Range("C4").Formula = "=COUNTA(Summary_1!A:A)-1"
$ws.Range("D4").Formula = "=B4/C4"
$ws.Range("D4").NumberFormat = "0.00%"

$ws.Range("A5").Value = "Percentage"
$ws.Range("B5").Formula = "=COUNTA(Summary_1!D:D)-1"
$ws.Range("C5").Formula = "=COUNTA(Summary_1!A:A)-1"
$ws.Range("D5").Formula = "=B5/C5"
$ws.Range("D5").NumberFormat = "0.00%"

# Row 6 has no label in column A.
$ws.Range("B6").Formula = "=COUNTA(Summary_1!E:E)-1"
$ws.Range("C6").Formula = "=COUNTA(Summary_1!A:A)-1"
$ws.Range("D6").Formula = "=B6/C6"
$ws.Range("D6").NumberFormat = "0.00%"
